# Performance doc update: refreshed benchmark numbers on the "Concise" sheet
# (Lines/Bytes/Whitespace/Alphanumeric/Special columns), a new "Comments"
# column (H), the resulting ratio-table recalculation, a rescroll/rezoom of
# the sheet view, and a small reposition of the second ("ratio") chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concise")

# --- Row 3 : Python -------------------------------------------------------
$ws.Range("C3").Value = 222
$ws.Range("D3").Value = 5408
$ws.Range("E3").Value = 2086
$ws.Range("F3").Value = 2595
$ws.Range("G3").Value = 679
$ws.Range("H3").Value = 50

# --- Row 4 : Bau -----------------------------------------------------------
$ws.Range("C4").Value = 227
$ws.Range("D4").Value = 5416
$ws.Range("E4").Value = 2063
$ws.Range("F4").Value = 2716
$ws.Range("G4").Value = 637
$ws.Range("H4").Value = 0

# --- Row 5 : Swift -----------------------------------------------------------
$ws.Range("D5").Value = 6566
$ws.Range("E5").Value = 2542
$ws.Range("G5").Value = 833
$ws.Range("H5").Value = 0

# --- Row 6 : Kotlin ----------------------------------------------------------
$ws.Range("H6").Value = 0

# --- Row 7 : C -----------------------------------------------------------
$ws.Range("D7").Value = 7268
$ws.Range("E7").Value = 2604
$ws.Range("F7").Value = 3477
$ws.Range("H7").Value = 60

# --- Row 8 : Go -----------------------------------------------------------
$ws.Range("D8").Value = 7037
$ws.Range("E8").Value = 2764
$ws.Range("F8").Value = 3217
$ws.Range("G8").Value = 1020
$ws.Range("H8").Value = 9

# --- Row 9 : Java -----------------------------------------------------------
$ws.Range("G9").Value = 1007
$ws.Range("H9").Value = 0

# --- Row 10 : Rust -----------------------------------------------------------
$ws.Range("D10").Value = 8311
$ws.Range("E10").Value = 3311
$ws.Range("F10").Value = 3693
$ws.Range("G10").Value = 1245
$ws.Range("H10").Value = 65

# --- Sheet view: rescroll + rezoom + reselect ------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 181
$ws.Range("F11").Select()

# --- Reposition the ratio chart ("Chart 5") --------------------------------
$chartObj = $ws.ChartObjects(2)
$chartObj.Left = 334.1764763779528
$chartObj.Top = 176.00007874015748
$chartObj.Width = 312.6360236220472
$chartObj.Height = 231.50000000000003
